$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update backlog status for the "Documentação" requirement (row 6): ---
# it moved from "Em andamento" to "Concluído", so recolor it green (same
# styling already used by the other "Concluído" cells, e.g. D8) and update
# the text.
$ws.Range("D8").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D6").Value = "Concluído"

# --- Update backlog status for the "Contexto" requirement (row 20): ---
# it moved from "Pendente" to "Concluído" as well.
$ws.Range("D8").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = "Concluído"

$excel.CutCopyMode = 0

# --- Remove the stray helper cell F9 that isn't part of the table anymore ---
$ws.Range("F9").Clear()

# --- Update the saved view/selection state ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("E22").Select()
